$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.670.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.124.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.34%  "
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5286"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.07"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09100"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.182"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.135.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.854"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.094"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001177"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.013"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06716"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.29%  "
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.340"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.731.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.35%  "
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.375.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.566"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.199"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1080"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.658"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.371"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.024"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.134"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02659"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06896"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2321"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6946"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.283"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.14%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.336"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6475"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000367"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.23%  "
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07316"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.46%  "
